# Regenerate save_data to use K (strikeouts) instead of Strike# column values.
# This updates column G ("K") for rows 2-18 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 2
    18 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
